$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new snapshot columns before the existing "B" (Jun_13) / "C" (Jun_10)
# columns. This shifts the old B -> D and old C -> E (values, formatting, and
# styles travel with the cells), leaving B:C empty for the two new snapshots.
$ws.Columns("B:C").Insert()

# Match the narrow width used by the rest of the snapshot columns (~8 chars).
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14

# New column headers (row 1): newest snapshot goes in B, next in C.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two new snapshot columns with the "UN" (unchanged) placeholder used
# throughout the sheet for rows with no rating action on that date.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
